# Remove the "CAVESITE_LOG_DIR=/var/www/logs" line from the dev
# installation instructions. The paragraph's own text plus its
# trailing paragraph mark are deleted, which merges it with the
# (empty) paragraph that precedes it.

$d = $word.ActiveDocument

$hit = $d.Content
$found = $hit.Find.Execute("CAVESITE_LOG_DIR=/var/www/logs")

if ($found) {
    # Extend the found range by one character to swallow the
    # paragraph mark that ends this paragraph, then delete the lot.
    $toDelete = $d.Range($hit.Start, $hit.End + 1)
    $toDelete.Delete()
}
